# Add a new roster row (Bryan Boh / b.boh@digipen.edu) to the team sheet,
# including an e-mail hyperlink on the new cell, matching the upload diff.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New row 6 data
$ws.Range("A6").Value = "Bryan Boh"
$ws.Range("B6").Value = "b.boh@digipen.edu"

# Hyperlink the new e-mail address, same pattern as the other email cells
$ws.Hyperlinks.Add($ws.Range("B6"), "mailto:b.boh@digipen.edu") | Out-Null

# Move the active selection to B7, where the author's cursor ended up next
$ws.Range("B7").Select() | Out-Null
